$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A (shifts B..L to C..M), so the
# table gets a new leading "No. Nota" column.
$ws.Columns("A").Insert()

# Widen the "Nama Pelanggan" column (now column C) a bit to fit its header.
$ws.Columns("C").ColumnWidth = 15.5

# New header cell for the inserted column.
$ws.Range("A7").Value = "No. Nota"

# Give A7 the same thin border used by the rest of the header row, but
# without the centered alignment the other header cells have.
$ws.Range("A7").Borders.LineStyle = 1

# Fix the "Kode Bbarang" typo -> "Kd. barang" (now column D after the insert).
$ws.Range("D7").Value = "Kd. barang"

# Match the saved selection state from the authored workbook.
$ws.Range("D7").Select()
